$d = $word.ActiveDocument

# "Versi" + "on" -> merge into a single run reading "Version"
# (limited to exactly this span so the surrounding proofErr spellStart/
# spellEnd markers are left untouched).
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version", 2)

# " 2" -> " 1." ; stop right before the _GoBack bookmark so the
# bookmarkStart/bookmarkEnd pair is not swallowed by the replace.
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 1.", 2)

# The trailing "." run (now redundant, after the bookmark) is deleted.
$r = $d.Range(10, 11)
$r.Delete()
